# Update the "Phones" (A) and "Prices" (B) columns for rows 2-25 with refreshed
# scrape data, mirroring the upstream test fixture regeneration.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values are plain text (phone names) - assign directly.
$ws.Range("A2").Value = 'Apple iPhone 12 (256GB) - White'
$ws.Range("A3").Value = 'Apple iPhone 13 (128GB) - Blue'
$ws.Range("A4").Value = 'Redmi A1 (Light Blue, 2GB RAM, 32GB Storage) | Segment Best AI Dual Cam | 5000mAh Battery | Leather Texture Design | Android 12'
$ws.Range("A5").Value = 'Redmi 10A (Slate Grey, 4GB RAM, 64GB Storage) | 2 Ghz Octa Core Helio G25 | 5000 mAh Battery | Finger Print Sensor | Upto 5GB RAM with RAM Booster'
$ws.Range("A6").Value = 'OPPO A74 5G (Fantastic Purple,6GB RAM,128GB Storage) with No Cost EMI/Additional Exchange Offers'
$ws.Range("A7").Value = 'Samsung Galaxy M04 Dark Blue, 4GB RAM, 64GB Storage | Upto 8GB RAM with RAM Plus | MediaTek Helio P35 | 5000 mAh Battery'
$ws.Range("A8").Value = 'Redmi A1 (Black, 2GB RAM, 32GB Storage) | Segment Best AI Dual Cam | 5000mAh Battery | Leather Texture Design | Android 12'
$ws.Range("A9").Value = 'OPPO A31 (Mystery Black, 6GB RAM, 128GB Storage) with No Cost EMI/Additional Exchange Offers'
$ws.Range("A10").Value = 'Redmi 9A Sport (Coral Green, 2GB RAM, 32GB Storage) | 2GHz Octa-core Helio G25 Processor | 5000 mAh Battery'
$ws.Range("A11").Value = 'Samsung Galaxy M13 (Midnight Blue, 4GB, 64GB Storage) | 6000mAh Battery | Upto 8GB RAM with RAM Plus'
$ws.Range("A12").Value = '(Renewed) OnePlus 6T (Mirror Black, 6GB RAM, 128GB Storage)'
$ws.Range("A13").Value = 'Apple iPhone 14 Plus 128GB (Product) RED'
$ws.Range("A14").Value = 'Samsung Galaxy M13 5G (Aqua Green, 4GB, 64GB Storage) | 5000mAh Battery | Upto 8GB RAM with RAM Plus'
$ws.Range("A15").Value = 'Samsung Galaxy M13 (Aqua Green, 4GB, 64GB Storage) | 6000mAh Battery | Upto 8GB RAM with RAM Plus'
$ws.Range("A16").Value = 'Redmi 10 Power (Sporty Orange, 8GB RAM, 128GB Storage)'
$ws.Range("A17").Value = 'Redmi 10A Sport (Sea Blue, 6GB RAM, 128GB Storage) | 2 Ghz Octa Cor Helio G25 | 5000 mAh Battery | Finger Print Sensor | Upto 8GB RAM with RAM Booster'
$ws.Range("A18").Value = 'OPPO A74 5G (Fluid Black, 6GB RAM, 128GB Storage) with No Cost EMI/Additional Exchange Offers'
$ws.Range("A19").Value = 'Samsung Galaxy M04 Light Green, 4GB RAM, 64GB Storage | Upto 8GB RAM with RAM Plus | MediaTek Helio P35 | 5000 mAh Battery'
$ws.Range("A20").Value = 'OnePlus Nord CE 2 Lite 5G (Black Dusk, 6GB RAM, 128GB Storage)'
$ws.Range("A21").Value = 'Redmi 11 Prime 5G (Chrome Silver, 4GB RAM 64GB ROM) | Prime Design | MTK Dimensity 700 | 50 MP Dual Cam | 5000mAh | 7 Band 5G'
$ws.Range("A22").Value = 'Apple iPhone 14 Pro 128GB Deep Purple'
$ws.Range("A23").Value = 'Tecno Spark 9 (Infinity Black, 6GB RAM,128GB Storage) | 11GB Expandable RAM | Helio G37 Gaming Processor'
$ws.Range("A24").Value = 'Tecno Spark 9 (Infinity Black, 6GB RAM,128GB Storage) | 11GB Expandable RAM | Helio G37 Gaming Processor'
$ws.Range("A25").Value = 'SOOPII "Grade-A Quality DM01 Pro Dual Support Desktop Mobile Phone + Tablet Tabletop Stand, Mobile Holder, Adjustable & Foldable, Aluminum Holder for Mobile Phone and Tablet (Up to 12 inch)'

# Column B values are digit/comma strings ("17,499") that Excel would
# otherwise auto-convert to numbers if assigned directly via .Value. Route them
# through a text formula first, then Copy/PasteSpecial(values) back onto
# themselves so the final stored cell is a literal shared-string, matching the
# original file's cell typing (t="s") instead of becoming t="n" numbers.
$ws.Range("B2").Formula = '="67,039"'
$ws.Range("B3").Formula = '="60,900"'
$ws.Range("B4").Formula = '="6,499"'
$ws.Range("B5").Formula = '="8,999"'
$ws.Range("B6").Formula = '="15,490"'
$ws.Range("B7").Formula = '="60,900"'
$ws.Range("B8").Formula = '="59,900"'
$ws.Range("B9").Formula = '="1,22,999"'
$ws.Range("B10").Formula = '="74,900"'
$ws.Range("B12").Formula = '="8,999"'
$ws.Range("B13").Formula = '="6,499"'
$ws.Range("B14").Formula = '="12,490"'
$ws.Range("B15").Formula = '="6,499"'
$ws.Range("B16").Formula = '="10,999"'
$ws.Range("B17").Formula = '="15,999"'
$ws.Range("B18").Formula = '="79,999"'
$ws.Range("B19").Formula = '="11,999"'
$ws.Range("B20").Formula = '="34,999"'
$ws.Range("B21").Formula = '="9,699"'
$ws.Range("B22").Formula = '="11,999"'
$ws.Range("B23").Formula = '="10,499"'
$ws.Range("B24").Formula = '="209"'
$ws.Range("B25").Formula = '="279"'

$ws.Range("B2:B25").Copy()
$ws.Range("B2:B25").PasteSpecial(-4163)

# B11 becomes an empty string in the target data.
$ws.Range("B11").Value = ""
